$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Add the new note text in B33, next to the existing "git diff" note in A33
$ws.Range("B33").Value = "what difference?"

# Update selection / scroll position to match the authored view
$ws.Range("B33").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

# Adjust the workbook window width recorded in the saved view
$excel.ActiveWindow.Width = 20730
